$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Document Type")

$ws.Range("A62").Value = "SG PEPPOL BIS Billing 3.0 Invoice"
$ws.Range("B62").Value = "busdox-docid-qns"
$ws.Range("C62").Value = "urn:oasis:names:specification:ubl:schema:xsd:Invoice-2::Invoice##urn:cen.eu:en16931:2017#conformant#urn:fdc:peppol.eu:2017:poacc:billing:international:sg:3.0::2.1"
$ws.Range("D62").Value = 4
$ws.Range("E62").Formula = "=FALSE"

$ws.Range("A63").Value = "SG PEPPOL BIS Billing 3.0 Credit Note"
$ws.Range("B63").Value = "busdox-docid-qns"
$ws.Range("C63").Value = "urn:oasis:names:specification:ubl:schema:xsd:Invoice-2::CreditNote##urn:cen.eu:en16931:2017#conformant#urn:fdc:peppol.eu:2017:poacc:billing:international:sg:3.0::2.1"
$ws.Range("D63").Value = 4
$ws.Range("E63").Formula = "=FALSE"
